# custom accuracy + 데이터 1000개
# - Round the numeric sensor readings in row 5 to 2 decimal places
# - Remove row 6 entirely (dimension shrinks from A1:AH6 to A1:AH5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to their 2-decimal-place rounded equivalents.
# (A5 is a date/time value and is left untouched; N5, S5, Z5, AB5 were
# already at 2 decimals so their values are unchanged.)
$ws.Range("B5").Value = 12.01
$ws.Range("C5").Value = 8.69
$ws.Range("D5").Value = 0.96
$ws.Range("E5").Value = 25.82
$ws.Range("F5").Value = 21.38
$ws.Range("G5").Value = 9.449999999999999
$ws.Range("H5").Value = 38.86
$ws.Range("I5").Value = 14.54
$ws.Range("J5").Value = 6.38
$ws.Range("K5").Value = 9.59
$ws.Range("L5").Value = 10.41
$ws.Range("M5").Value = 10.86
$ws.Range("N5").Value = 3.02
$ws.Range("O5").Value = 9.4
$ws.Range("P5").Value = 13.33
$ws.Range("Q5").Value = 8.01
$ws.Range("R5").Value = 0.8
$ws.Range("S5").Value = 0.58
$ws.Range("T5").Value = 135.75
$ws.Range("U5").Value = 26.35
$ws.Range("V5").Value = 8.68
$ws.Range("W5").Value = 17.59
$ws.Range("X5").Value = 9.41
$ws.Range("Y5").Value = 1.15
$ws.Range("Z5").Value = 18.37
$ws.Range("AA5").Value = 7.66
$ws.Range("AB5").Value = 6.87
$ws.Range("AC5").Value = 8.06
$ws.Range("AD5").Value = 10.9
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 35.13
$ws.Range("AG5").Value = 4.84
$ws.Range("AH5").Value = 10.85

# Delete row 6 entirely - this shifts the dimension ref down to A1:AH5
$ws.Rows(6).Delete()
